# Apply updated gene-expression values for the fresh-install rerun (2025-05-08).
# The underlying R analysis produced slightly different counts and reordered
# several genes by Mean_GeneSymbol_TPM, and added new rows for genes that
# newly exceeded the reporting threshold.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Tff1", 96396, $true, $true, $true, $true),
    @("mt-Rnr1", 94140, $true, $false, $false, $false),
    @("mt-Rnr2", 60108, $true, $false, $false, $false),
    @("Yam1", 54893, $true, $false, $false, $false),
    @("Gkn2", 37154, $true, $false, $true, $false),
    @("Gkn1", 36085, $true, $false, $true, $false),
    @("Psca", 17428, $true, $false, $true, $false),
    @("mt-Atp6", 13473, $true, $false, $false, $false),
    @("Dpcr1", 9829, $true, $false, $false, $false),
    @("Krt8", 6479, $true, $false, $false, $false),
    @("mt-Co1", 6061, $true, $false, $false, $false),
    @("Gm26917", 5763, $true, $false, $false, $false),
    @("Muc5ac", 4931, $true, $true, $true, $true),
    @("Lgals2", 4805, $true, $false, $true, $false),
    @("Krt19", 4695, $true, $false, $false, $false),
    @("Rn7sk", 4561, $true, $false, $false, $false),
    @("Car2", 4260, $true, $false, $false, $false),
    @("Gm23935", 4228, $true, $false, $false, $false),
    @("mt-Cytb", 4202, $true, $false, $false, $false),
    @("Gsta4", 4158, $true, $false, $false, $false),
    @("Gsta1", 3992, $true, $false, $false, $false),
    @("Gsto1", 3949, $true, $false, $false, $false),
    @("mt-Co2", 3829, $true, $false, $false, $false),
    @("Ahnak", 3771, $true, $false, $false, $false),
    @("Sprr2a2", 3728, $true, $false, $false, $false),
    @("Gsdma2", 1982, $false, $false, $true, $false),
    @("Pgc", 1404, $false, $false, $true, $false),
    @("Anxa10", 961, $false, $false, $true, $false),
    @("Tff2", 947, $false, $false, $true, $false),
    @("Ptgr1", 800, $false, $false, $true, $false),
    @("Sprr2a1", 754, $false, $false, $true, $false),
    @("Krt20", 665, $false, $true, $true, $false),
    @("Lgals4", 590, $false, $true, $true, $false),
    @("Muc1", 372, $false, $false, $true, $false),
    @("Cldn18", 301, $false, $true, $true, $false),
    @("Sult1c2", 263, $false, $false, $true, $false),
    @("Clps", 254, $false, $false, $true, $false),
    @("Sprr2a3", 189, $false, $false, $true, $false),
    @("Foxq1", 161, $false, $false, $true, $false),
    @("Pla2g1b", 119, $false, $false, $true, $false),
    @("Il33", 43, $false, $true, $true, $true),
    @("Akr1b10", 9, $false, $true, $true, $false)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $r = $data[$i]
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $ws.Cells.Item($row, 5).Value = $r[4]
    $ws.Cells.Item($row, 6).Value = $r[5]
}

$lastRow = $startRow + $data.Count - 1
Write-Output "Updated $($data.Count) data rows; sheet now spans A1:F$lastRow"
